# Release Checklist v2.4 update:
# Mark several more checklist rows as "Done" instead of "Open", and move the
# current on-screen selection down to C31 (scrolled so row 2 is back at the
# top of the view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Status column ("Open" -> "Done") for the newly-completed checklist items
$ws.Range("C14").Value = "Done"
$ws.Range("C20").Value = "Done"
$ws.Range("C24").Value = "Done"
$ws.Range("C25").Value = "Done"
$ws.Range("C26").Value = "Done"
$ws.Range("C27").Value = "Done"
$ws.Range("C30").Value = "Done"

# Restore the view to the top of the sheet and move the active selection
# down to the (now blank) row just past the last checklist item.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C31").Select()
